# Trade #12 closed/logged: leadlag DOWN trade opened 2026-02-16 21:54:45.
# Appends the new trade row to both the "All Trades" ledger sheet and the
# strategy-specific "leadlag" sheet (mirrors how prior trades are duplicated
# across both sheets in this workbook).

$wb = $excel.ActiveWorkbook

function Add-TradeRow {
    param($Sheet, $Row)

    # Column A - Trade #
    $Sheet.Cells.Item($Row, 1).Value = 12

    # Column B - Date (text that looks like a date -> force Text format so
    # Excel doesn't silently convert it to a date serial number, then strip
    # the format back off so the cell keeps the workbook's default style).
    $cell = $Sheet.Cells.Item($Row, 2)
    $cell.NumberFormat = "@"
    $cell.Value = "2026-02-16"
    $cell.ClearFormats()

    # Column C - Time (same text-vs-autoconvert concern as the date column).
    $cell = $Sheet.Cells.Item($Row, 3)
    $cell.NumberFormat = "@"
    $cell.Value = "21:54:45"
    $cell.ClearFormats()

    # Column D - Strategy
    $Sheet.Cells.Item($Row, 4).Value = "leadlag"

    # Column E - Side
    $Sheet.Cells.Item($Row, 5).Value = "DOWN"

    # Column F - Entry Price
    $Sheet.Cells.Item($Row, 6).Value = 68238.355

    # Column G - Exit Price (still open -> blank/empty)
    $Sheet.Cells.Item($Row, 7).Value = ""

    # Column H - Status
    $Sheet.Cells.Item($Row, 8).Value = "OPEN"

    # Column I - P&L %
    $Sheet.Cells.Item($Row, 9).Value = 0

    # Column J - P&L $
    $Sheet.Cells.Item($Row, 10).Value = 0

    # Column K - Capital After
    $Sheet.Cells.Item($Row, 11).Value = 100

    # Column L - Confidence
    $Sheet.Cells.Item($Row, 12).Value = 0.75

    # Column M - Entry Reason
    $Sheet.Cells.Item($Row, 13).Value = "Coinbase leading with -0.124% move"

    # Column N - Exit Reason (still open -> blank/empty)
    $Sheet.Cells.Item($Row, 14).Value = ""

    # Column O - Duration (min)
    $Sheet.Cells.Item($Row, 15).Value = 0
}

# "All Trades" combined ledger: new trade becomes row 13 (previous last row
# was 12).
$allTrades = $wb.Worksheets.Item("All Trades")
Add-TradeRow $allTrades 13

# "leadlag" strategy-specific sheet: new trade becomes row 12 (previous last
# row was 11).
$leadlag = $wb.Worksheets.Item("leadlag")
Add-TradeRow $leadlag 12
